$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the stale "1618 lake sims parkway oc" sheet entirely.
$wb.Sheets.Item("1618 lake sims parkway oc").Delete() | Out-Null

# Update the content on the remaining address sheet ("451 clear blue way mcdono")
# to reflect the refreshed sex-offender registry figures and blank out the
# now-unavailable Bestplaces crime stats.
$ws = $wb.Sheets.Item("451 clear blue way mcdono")
$ws.Range("B10").Value = "According to our research of Georgia and other state lists, there were 97 registered sex offenders living in McDonough, Georgia as of December 22, 2020.`nThe ratio of all residents to sex offenders in McDonough is 243 to 1."
# The multi-line text would otherwise auto-expand row 10's height; put it back
# to the sheet's normal (non-custom) row height, matching the source file.
$ws.Rows.Item(10).AutoFit() | Out-Null
$ws.Range("B23").Value = "NA"
$ws.Range("B24").Value = "NA"
$ws.Range("B25").Value = "NA"

# Make "test" the active/selected sheet (was "1618 lake sims parkway oc").
$wb.Sheets.Item("test").Activate() | Out-Null
